# Stage 1: update companies data
#
# Applies the row re-shuffle described in the diff: the A/B/H/I/J/K
# columns for several data rows are replaced with another row's values
# (companies were re-sorted), while columns C-G (dates/status/source/
# time) are left untouched since they are identical between the swapped
# rows.
#
# Every cell in this sheet is stored as text (inlineStr), including the
# numeric-looking "Company Number" column. Excel's COM `.Value` setter
# auto-coerces digit-only / date-looking strings into real numbers or
# dates, which would change the cell's underlying type. To avoid that we
# temporarily force the cell's number format to Text ("@") before the
# assignment, then restore the original Style object afterwards so no
# visible formatting/style attribute is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Addr, $Text) {
    $cell = $ws.Range($Addr)
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = $originalStyle
}

# row -> hashtable of column letter -> new text value
$edits = [ordered]@{
    3 = [ordered]@{
        A = "GANDER INVESTMENTS LTD"
        B = "16473515"
        H = "Investments"
        I = "68100,68209"
        J = ""
        K = ""
    }
    4 = [ordered]@{
        A = "SEVEN (HOLDCO) LIMITED"
        B = "16473606"
        H = "Other"
        I = "64209"
        J = "Activities of other holding companies n.e.c."
        K = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles."
    }
    5 = [ordered]@{
        A = "AJ INVESTMENT AND CONSULTANCY LTD"
        B = "16473328"
        H = "Investments"
        I = "64306,70229"
        J = "Activities of real estate investment trusts"
        K = "UK-regulated REIT companies."
    }
    6 = [ordered]@{
        A = "BRIDGEWICK PARTNERS LIMITED"
        B = "16473142"
        H = "Partners"
        I = "64999"
        J = "Financial intermediation not elsewhere classified"
        K = "Catch-all credit-oriented SPVs for novel lending structures."
    }
    7 = [ordered]@{
        A = "MARMIMI HOLDING LIMITED"
        B = "16473234"
        H = "Other"
    }
    9 = [ordered]@{
        A = "TLJ INVESTMENT LTD"
        B = "16473151"
        I = "41100,55100,68100"
        J = ""
        K = ""
    }
    10 = [ordered]@{
        A = "INTERCONTINENTAL HOLDING COMPANY LIMITED"
        B = "16473418"
    }
    11 = [ordered]@{
        A = "GAUNT CAPITAL LTD"
        B = "16473262"
        H = "Capital"
        I = "64209"
        J = "Activities of other holding companies n.e.c."
        K = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles."
    }
}

foreach ($rowNum in $edits.Keys) {
    $cols = $edits[$rowNum]
    foreach ($col in $cols.Keys) {
        $addr = "$col$rowNum"
        Set-TextValue $addr $cols[$col]
    }
}
